$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in F2 per PR feedback
$ws.Range("F2").Value = 7

# Update the active selection to F2 (was F3)
$ws.Range("F2").Select()
